# Commit: "27.06.19 Today sales details"
# Update the "Raju Ahamed" sheet: change the date from 26.06.19 to 27.06.19
# and update the today's sales quantities (and their dependent formulas
# recompute automatically), then move the saved view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raju Ahamed")

# --- Update the date label (shared between the two copies of the sheet,
# both A4 and A31 show "Date: 26.06.19" -> "Date: 27.06.19") ---
$ws.Range("A4").Value = "Date: 27.06.19"
$ws.Range("A31").Value = "Date: 27.06.19"

# --- Update the "Qty." column for the first table (rows 6-11) ---
$ws.Range("E6").Value = 34
$ws.Range("E7").Value = 122
$ws.Range("E8").Value = 509
$ws.Range("E9").Value = 54
$ws.Range("E10").Value = 70
$ws.Range("E11").Value = 100

# --- Update the "Qty." column for the second (duplicate) table (rows 33-38) ---
$ws.Range("E33").Value = 34
$ws.Range("E34").Value = 122
$ws.Range("E35").Value = 509
$ws.Range("E36").Value = 54
$ws.Range("E37").Value = 70
$ws.Range("E38").Value = 100

# --- Move the saved scroll position / selection on the sheet ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.SetTopLeftVisibleCell("A29")
$ws.Range("A32").Select()
